$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume change (E) columns for rows 2-50.
# D-column values are numeric-looking text (e.g. "1.01", "25.786.61") that must
# remain stored as text, matching the original inline-string cell type. We
# temporarily force a text number format before assigning the value, then
# restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.786.61'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -0.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.586.67'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -1.92%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.27'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -1.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.01'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.482'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -3.49%  '

$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0616'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.08'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -1.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.808.63'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -1.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.587.46'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -2.04%  '

$ws.Range("E14").Value = '  -2.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.510'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -2.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.804.40'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -0.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0721'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -1.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '59.92'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -2.44%  '

$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.19'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.18'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -1.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.37'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -1.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.92'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -1.62%  '

$ws.Range("E24").Value = '  -1.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.61'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.69'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -1.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.12'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.45'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -2.74%  '

$ws.Range("E30").Value = '  -5.53%  '

$ws.Range("E31").Value = '  -1.13%  '

$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.02'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -2.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.49'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +0.46%  '

$ws.Range("E35").Value = '  -1.89%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.102.07'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -2.54%  '

$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.35'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -1.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.504'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -1.33%  '

$ws.Range("E40").Value = '  -1.84%  '

$ws.Range("E41").Value = '  -7.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.811'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +8.24%  '

$ws.Range("E43").Value = '  +2.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.80'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -4.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.723.32'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -1.73%  '

$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.50'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.34'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -1.19%  '

$ws.Range("E49").Value = '  -1.78%  '

$ws.Range("E50").Value = '  -0.52%  '
